$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.779.97"
$ws.Range("E2").Value = "  +1.22%  "

# Row 3
$ws.Range("D3").Value = "3.408.66"
$ws.Range("E3").Value = "  +0.63%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.28%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "562.85"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.25%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "176.97"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.27%  "

# Row 7
$ws.Range("E7").Value = "  +0.34%  "

# Row 8
$ws.Range("D8").Value = "3.402.07"
$ws.Range("E8").Value = "  +0.52%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.00"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.03%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.172"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +3.74%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.643"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.11%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.78"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.76%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000279"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.47%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.07%  "

# Row 15
$ws.Range("D15").Value = "3.948.61"
$ws.Range("E15").Value = "  +1.16%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.38"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.03%  "

# Row 17
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "3.430.49"
$ws.Range("E17").Value = "  +2.20%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.120"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.14%  "

# Row 19
$ws.Range("D19").Value = "65.673.17"
$ws.Range("E19").Value = "  +1.55%  "

# Row 20
$ws.Range("E20").Value = "  -0.11%  "

# Row 21
$ws.Range("E21").Value = "  +1.41%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "490.45"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +4.30%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.95"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.67%  "

# Row 24
$ws.Range("E24").Value = "  -0.09%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "89.61"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.39%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "14.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.77%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.06%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.40%  "

# Row 29
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "31.64"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.69%  "

# Row 30
$ws.Range("B30").Value = "Filecoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.77"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.37%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.58"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "11.54"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.14%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "62.52"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "576.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.60%  "

# Row 35
$ws.Range("E35").Value = "  -0.57%  "

# Row 36
$ws.Range("E36").Value = "  -0.09%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.61"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.20%  "

# Row 38
$ws.Range("E38").Value = "  +0.34%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.05%  "

# Row 40
$ws.Range("E40").Value = "  +0.54%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0742"
$ws.Range("E41").Value = "  -2.81%  "

# Row 42
$ws.Range("D42").Value = "3.126.11"
$ws.Range("E42").Value = "  +0.91%  "

# Row 43
$ws.Range("E43").Value = "  -1.03%  "

# Row 44
$ws.Range("E44").Value = "  +1.33%  "

# Row 45
$ws.Range("E45").Value = "  +1.71%  "

# Row 46
$ws.Range("B46").Value = "Fetch.AI"
$ws.Range("C46").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.45"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.36%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.16"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.58%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "140.49"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.57"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.15%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.66%  "
